$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 196.36363
$ws.Range("I9").Value = 229.16667
$ws.Range("J9").Value = 157
$ws.Range("K9").Value = 229.16667
$ws.Range("L9").Value = 157
$ws.Range("M9").Value = -60.16667000000001
$ws.Range("N9").Value = -495

$ws.Range("H15").Value = 872.1739
$ws.Range("I15").Value = 872.1739
$ws.Range("K15").Value = 2616.5217
$ws.Range("M15").Value = -2447.5217

$ws.Range("H40").Value = 3155.5293
$ws.Range("I40").Value = 2682.625
$ws.Range("J40").Value = 3575.889
$ws.Range("K40").Value = 2682.625
$ws.Range("L40").Value = 3575.889
$ws.Range("M40").Value = -2507.625
$ws.Range("N40").Value = -3925.889

$ws.Range("H70").Value = 2817.5454
$ws.Range("I70").Value = 2599.1667
$ws.Range("J70").Value = 3079.6
$ws.Range("K70").Value = 7797.500100000001
$ws.Range("L70").Value = 9238.799999999999
$ws.Range("M70").Value = -7527.500100000001
$ws.Range("N70").Value = -9778.799999999999

$ws.Range("H73").Value = 2817.5454
$ws.Range("I73").Value = 2599.1667
$ws.Range("J73").Value = 3079.6
$ws.Range("K73").Value = 7797.500100000001
$ws.Range("L73").Value = 9238.799999999999
$ws.Range("M73").Value = -6861.500100000001
$ws.Range("N73").Value = -11110.8

$ws.Range("H80").Value = 658.3333
$ws.Range("I80").Value = 637.5
$ws.Range("J80").Value = 700
$ws.Range("K80").Value = 1912.5
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -914.5
$ws.Range("N80").Value = -4096

$ws.Range("H83").Value = 658.3333
$ws.Range("I83").Value = 637.5
$ws.Range("J83").Value = 700
$ws.Range("K83").Value = 5737.5
$ws.Range("L83").Value = 6300
$ws.Range("M83").Value = -745.5
$ws.Range("N83").Value = -16284

$ws.Range("H107").Value = 622.1111
$ws.Range("I107").Value = 654.2941
$ws.Range("J107").Value = 75
$ws.Range("K107").Value = 654.2941
$ws.Range("L107").Value = 75
$ws.Range("M107").Value = 1265.7059
$ws.Range("N107").Value = -3915

$ws.Range("H116").Value = 4999.75
$ws.Range("J116").Value = 4999.75
$ws.Range("L116").Value = 4999.75
$ws.Range("N116").Value = -11883.75

$ws.Range("H125").Value = 17605.908
$ws.Range("I125").Value = 28420.54
$ws.Range("J125").Value = 1984.7778
$ws.Range("K125").Value = 255784.86
$ws.Range("L125").Value = 17863.0002
$ws.Range("M125").Value = -253324.86
$ws.Range("N125").Value = -22783.0002

$ws.Range("H132").Value = 2925.8125
$ws.Range("I132").Value = 2920.9333
$ws.Range("K132").Value = 8762.7999
$ws.Range("M132").Value = -6232.7999

$ws.Range("H135").Value = 464.12122
$ws.Range("I135").Value = 310.5484
$ws.Range("J135").Value = 2844.5
$ws.Range("K135").Value = 2794.9356
$ws.Range("L135").Value = 25600.5
$ws.Range("M135").Value = -259.9356000000002
$ws.Range("N135").Value = -30670.5

$ws.Range("H137").Value = 4223.5
$ws.Range("I137").Value = 3580.0232
$ws.Range("J137").Value = 5541.095
$ws.Range("K137").Value = 10740.0696
$ws.Range("L137").Value = 16623.285
$ws.Range("M137").Value = -8190.069600000001
$ws.Range("N137").Value = -21723.285

$ws.Range("H138").Value = 1675.9574
$ws.Range("I138").Value = 1060.742
$ws.Range("J138").Value = 2867.9375
$ws.Range("K138").Value = 3182.226
$ws.Range("L138").Value = 8603.8125
$ws.Range("M138").Value = 1957.774
$ws.Range("N138").Value = -18883.8125

$ws.Range("H141").Value = 19658.408
$ws.Range("J141").Value = 7498.778
$ws.Range("L141").Value = 22496.334
$ws.Range("N141").Value = -32856.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6123.75
$ws.Range("I2").Value = 6166.6665
$ws.Range("K2").Value = 6166.6665
$ws.Range("M2").Value = -6053.6665

$ws.Range("H5").Value = 281.5
$ws.Range("J5").Value = 283.8
$ws.Range("L5").Value = 283.8
$ws.Range("N5").Value = -507.8

$ws.Range("H32").Value = 17201.438
$ws.Range("I32").Value = 3239.7722
$ws.Range("J32").Value = 127498.6
$ws.Range("K32").Value = 3239.7722
$ws.Range("L32").Value = 127498.6
$ws.Range("M32").Value = -2952.7722
$ws.Range("N32").Value = -128072.6

$ws.Range("H45").Value = 534018.6
$ws.Range("I45").Value = 843435.8
$ws.Range("J45").Value = 3589.1428
$ws.Range("K45").Value = 843435.8
$ws.Range("L45").Value = 3589.1428
$ws.Range("M45").Value = -843058.8
$ws.Range("N45").Value = -4343.1428

$ws.Range("H61").Value = 1527.7858
$ws.Range("I61").Value = 1447.7273
$ws.Range("J61").Value = 1821.3334
$ws.Range("K61").Value = 1447.7273
$ws.Range("L61").Value = 1821.3334
$ws.Range("M61").Value = -1235.7273
$ws.Range("N61").Value = -2245.3334

$ws.Range("H74").Value = 2509.9167
$ws.Range("I74").Value = 1691.125
$ws.Range("K74").Value = 1691.125
$ws.Range("M74").Value = -817.125

$ws.Range("H77").Value = 2509.9167
$ws.Range("I77").Value = 1691.125
$ws.Range("K77").Value = 8455.625
$ws.Range("M77").Value = -4087.625

$ws.Range("H97").Value = 701.7222
$ws.Range("I97").Value = 684.17645
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 684.17645
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -188.17645
$ws.Range("N97").Value = -1992

$ws.Range("H102").Value = 1399
$ws.Range("I102").Value = 1399
$ws.Range("K102").Value = 1399
$ws.Range("M102").Value = 223

$ws.Range("H110").Value = 1265.0476
$ws.Range("I110").Value = 1122.4375
$ws.Range("K110").Value = 1122.4375
$ws.Range("M110").Value = 922.5625

$ws.Range("H116").Value = 6123.75
$ws.Range("I116").Value = 6166.6665
$ws.Range("K116").Value = 6166.6665
$ws.Range("M116").Value = -3872.6665

$ws.Range("H122").Value = 1606.8823
$ws.Range("I122").Value = 1539.5333
$ws.Range("J122").Value = 2112
$ws.Range("K122").Value = 4618.5999
$ws.Range("L122").Value = 6336
$ws.Range("M122").Value = -2168.5999
$ws.Range("N122").Value = -11236

$ws.Range("H131").Value = 79998.75
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 79998.75
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 79998.75
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -90078.75

$ws.Range("H132").Value = 1234.4615
$ws.Range("I132").Value = 1018.54285
$ws.Range("K132").Value = 3055.62855
$ws.Range("M132").Value = -525.6285500000004

$ws.Range("H136").Value = 1527.7858
$ws.Range("I136").Value = 1447.7273
$ws.Range("J136").Value = 1821.3334
$ws.Range("K136").Value = 4343.1819
$ws.Range("L136").Value = 5464.0002
$ws.Range("M136").Value = -1793.1819
$ws.Range("N136").Value = -10564.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6123.75
$ws.Range("I3").Value = 6166.6665
$ws.Range("K3").Value = 6166.6665
$ws.Range("M3").Value = -6052.6665

$ws.Range("H4").Value = 281.5
$ws.Range("J4").Value = 283.8
$ws.Range("L4").Value = 283.8
$ws.Range("N4").Value = -513.8

$ws.Range("H17").Value = 659.8333
$ws.Range("J17").Value = 659.8333
$ws.Range("L17").Value = 659.8333
$ws.Range("N17").Value = -1003.8333

$ws.Range("H20").Value = 5659.3076
$ws.Range("I20").Value = 4948.143
$ws.Range("J20").Value = 7469.5454
$ws.Range("K20").Value = 4948.143
$ws.Range("L20").Value = 7469.5454
$ws.Range("M20").Value = -4701.143
$ws.Range("N20").Value = -7963.5454

$ws.Range("H94").Value = 1227.027
$ws.Range("I94").Value = 1118.2122
$ws.Range("K94").Value = 1118.2122
$ws.Range("M94").Value = -667.2121999999999

$ws.Range("H99").Value = 2855.4375
$ws.Range("I99").Value = 3255
$ws.Range("K99").Value = 3255
$ws.Range("M99").Value = -1757

$ws.Range("H105").Value = 5026.3076
$ws.Range("I105").Value = 5321.1665
$ws.Range("K105").Value = 5321.1665
$ws.Range("M105").Value = -3574.1665

$ws.Range("H134").Value = 1055.4375
$ws.Range("I134").Value = 1055.4375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3166.3125
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -631.3125
$ws.Range("N134").ClearContents()

$ws.Range("H141").Value = 72304.89999999999
$ws.Range("J141").Value = 70203.875
$ws.Range("L141").Value = 70203.875
$ws.Range("N141").Value = -80563.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40469.383
$ws.Range("I31").Value = 41728.2
$ws.Range("J31").Value = 8999
$ws.Range("K31").Value = 41728.2
$ws.Range("L31").Value = 8999
$ws.Range("M31").Value = -41433.2
$ws.Range("N31").Value = -9589

$ws.Range("H33").Value = 8000
$ws.Range("I33").Value = 8000
$ws.Range("K33").Value = 8000
$ws.Range("M33").Value = -7621

$ws.Range("H34").Value = 40469.383
$ws.Range("I34").Value = 41728.2
$ws.Range("J34").Value = 8999
$ws.Range("K34").Value = 41728.2
$ws.Range("L34").Value = 8999
$ws.Range("M34").Value = -41526.2
$ws.Range("N34").Value = -9403

$ws.Range("H36").Value = 7749.8
$ws.Range("I36").Value = 7187.25
$ws.Range("K36").Value = 7187.25
$ws.Range("M36").Value = -6799.25

$ws.Range("H40").Value = 7749.8
$ws.Range("I40").Value = 7187.25
$ws.Range("K40").Value = 7187.25
$ws.Range("M40").Value = -7027.25

$ws.Range("H58").Value = 1977.4572
$ws.Range("I58").Value = 1710.125
$ws.Range("J58").Value = 2560.7273
$ws.Range("K58").Value = 1710.125
$ws.Range("L58").Value = 2560.7273
$ws.Range("M58").Value = -1507.125
$ws.Range("N58").Value = -2966.7273

$ws.Range("H62").Value = 6410.4287
$ws.Range("I62").Value = 6312.1665
$ws.Range("K62").Value = 6312.1665
$ws.Range("M62").Value = -5688.1665

$ws.Range("H65").Value = 6410.4287
$ws.Range("I65").Value = 6312.1665
$ws.Range("K65").Value = 31560.8325
$ws.Range("M65").Value = -28440.8325

$ws.Range("H122").Value = 7531.75
$ws.Range("I122").Value = 6854.4287
$ws.Range("K122").Value = 20563.2861
$ws.Range("M122").Value = -18113.2861

$ws.Range("H132").Value = 3997.195
$ws.Range("I132").Value = 3786
$ws.Range("K132").Value = 11358
$ws.Range("M132").Value = -8828

$ws.Range("H134").Value = 5223.8335
$ws.Range("I134").Value = 5244.273
$ws.Range("K134").Value = 15732.819
$ws.Range("M134").Value = -13197.819

$ws.Range("H136").Value = 1977.4572
$ws.Range("I136").Value = 1710.125
$ws.Range("J136").Value = 2560.7273
$ws.Range("K136").Value = 5130.375
$ws.Range("L136").Value = 7682.1819
$ws.Range("M136").Value = -2580.375
$ws.Range("N136").Value = -12782.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 90
$ws.Range("I7").Value = 90
$ws.Range("K7").Value = 270
$ws.Range("M7").Value = -158

$ws.Range("H34").Value = 92279.82000000001
$ws.Range("I34").Value = 142969.72
$ws.Range("J34").Value = 3572.5
$ws.Range("K34").Value = 428909.16
$ws.Range("L34").Value = 10717.5
$ws.Range("M34").Value = -428825.16
$ws.Range("N34").Value = -10885.5

$ws.Range("H60").Value = 200650.2
$ws.Range("J60").Value = 499.5
$ws.Range("L60").Value = 1498.5
$ws.Range("N60").Value = -2000.5

$ws.Range("H70").Value = 4229.636
$ws.Range("I70").Value = 3992.3333
$ws.Range("J70").Value = 4318.625
$ws.Range("K70").Value = 11976.9999
$ws.Range("L70").Value = 12955.875
$ws.Range("M70").Value = -11661.9999
$ws.Range("N70").Value = -13585.875

$ws.Range("H73").Value = 4229.636
$ws.Range("I73").Value = 3992.3333
$ws.Range("J73").Value = 4318.625
$ws.Range("K73").Value = 11976.9999
$ws.Range("L73").Value = 12955.875
$ws.Range("M73").Value = -10884.9999
$ws.Range("N73").Value = -15139.875

$ws.Range("H75").Value = 4118.625
$ws.Range("J75").Value = 5158.1665
$ws.Range("L75").Value = 15474.4995
$ws.Range("N75").Value = -17470.4995

$ws.Range("H78").Value = 4118.625
$ws.Range("J78").Value = 5158.1665
$ws.Range("L78").Value = 46423.4985
$ws.Range("N78").Value = -56407.4985

$ws.Range("H86").Value = 1161
$ws.Range("J86").Value = 1679.8
$ws.Range("L86").Value = 5039.4
$ws.Range("N86").Value = -7411.4

$ws.Range("H87").Value = 3664.6667
$ws.Range("I87").Value = 3664.6667
$ws.Range("K87").Value = 10994.0001
$ws.Range("M87").Value = -9746.000100000001

$ws.Range("H89").Value = 1161
$ws.Range("J89").Value = 1679.8
$ws.Range("L89").Value = 15118.2
$ws.Range("N89").Value = -26974.2

$ws.Range("H90").Value = 3664.6667
$ws.Range("I90").Value = 3664.6667
$ws.Range("K90").Value = 32982.0003
$ws.Range("M90").Value = -26742.0003

$ws.Range("H97").Value = 3185

$ws.Range("H107").Value = 663.3514
$ws.Range("J107").Value = 721.86664
$ws.Range("L107").Value = 2165.59992
$ws.Range("N107").Value = -6005.59992

$ws.Range("H109").Value = 94488.55
$ws.Range("I109").Value = 112152.664
$ws.Range("K109").Value = 336457.992
$ws.Range("M109").Value = -335417.992

$ws.Range("H131").Value = 131429.4
$ws.Range("I131").Value = 334103.34
$ws.Range("J131").Value = 44569.145
$ws.Range("K131").Value = 1002310.02
$ws.Range("L131").Value = 133707.435
$ws.Range("M131").Value = -997270.02
$ws.Range("N131").Value = -143787.435

$ws.Range("H132").Value = 1985.7778
$ws.Range("I132").Value = 1174.4
$ws.Range("K132").Value = 10569.6
$ws.Range("M132").Value = -8039.6

$ws.Range("H140").Value = 3607.625
$ws.Range("I140").Value = 2301.6924
$ws.Range("J140").Value = 9266.666999999999
$ws.Range("K140").Value = 6905.0772
$ws.Range("L140").Value = 27800.001
$ws.Range("M140").Value = -1725.0772
$ws.Range("N140").Value = -38160.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5242.222
$ws.Range("I70").Value = 5276.6
$ws.Range("K70").Value = 5276.6
$ws.Range("M70").Value = -5006.6

$ws.Range("H73").Value = 5242.222
$ws.Range("I73").Value = 5276.6
$ws.Range("K73").Value = 5276.6
$ws.Range("M73").Value = -4340.6

$ws.Range("H102").Value = 3591.2778
$ws.Range("I102").Value = 3840.1875
$ws.Range("K102").Value = 3840.1875
$ws.Range("M102").Value = -2218.1875

$ws.Range("H113").Value = 2829.6843
$ws.Range("I113").Value = 3107.5833
$ws.Range("J113").Value = 2353.2856
$ws.Range("K113").Value = 3107.5833
$ws.Range("L113").Value = 2353.2856
$ws.Range("M113").Value = -937.5832999999998
$ws.Range("N113").Value = -6693.2856

$ws.Range("H122").Value = 1679.6757
$ws.Range("I122").Value = 1571.8667
$ws.Range("J122").Value = 2141.7144
$ws.Range("K122").Value = 4715.6001
$ws.Range("L122").Value = 6425.1432
$ws.Range("M122").Value = -2265.6001
$ws.Range("N122").Value = -11325.1432

$ws.Range("H126").Value = 2844.3333
$ws.Range("I126").Value = 2488.2
$ws.Range("K126").Value = 7464.599999999999
$ws.Range("M126").Value = -4994.599999999999

$ws.Range("H132").Value = 16139.667
$ws.Range("I132").Value = 18524.389
$ws.Range("J132").Value = 1831.3334
$ws.Range("K132").Value = 55573.167
$ws.Range("L132").Value = 5494.0002
$ws.Range("M132").Value = -53043.167
$ws.Range("N132").Value = -10554.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 679
$ws.Range("I16").Value = 377.83334
$ws.Range("J16").Value = 1040.4
$ws.Range("K16").Value = 377.83334
$ws.Range("L16").Value = 1040.4
$ws.Range("M16").Value = -207.83334
$ws.Range("N16").Value = -1380.4

$ws.Range("H22").Value = 1243.9546
$ws.Range("I22").Value = 1148.9
$ws.Range("J22").Value = 1323.1666
$ws.Range("K22").Value = 1148.9
$ws.Range("L22").Value = 1323.1666
$ws.Range("M22").Value = -853.9000000000001
$ws.Range("N22").Value = -1913.1666

$ws.Range("H27").Value = 1243.9546
$ws.Range("I27").Value = 1148.9
$ws.Range("J27").Value = 1323.1666
$ws.Range("K27").Value = 1148.9
$ws.Range("L27").Value = 1323.1666
$ws.Range("M27").Value = -1041.9
$ws.Range("N27").Value = -1537.1666

$ws.Range("H55").Value = 193.69565
$ws.Range("I55").Value = 200.57143
$ws.Range("K55").Value = 200.57143
$ws.Range("M55").Value = -27.57142999999999

$ws.Range("H100").Value = 27113.8
$ws.Range("I100").Value = 5296.75
$ws.Range("K100").Value = 5296.75
$ws.Range("M100").Value = -4755.75

$ws.Range("H106").Value = 16104.857
$ws.Range("J106").Value = 16104.857
$ws.Range("L106").Value = 16104.857
$ws.Range("N106").Value = -18628.857

$ws.Range("H132").Value = 2542.4614
$ws.Range("I132").Value = 2158.6072
$ws.Range("J132").Value = 4930.8887
$ws.Range("K132").Value = 6475.821599999999
$ws.Range("L132").Value = 14792.6661
$ws.Range("M132").Value = -3945.821599999999
$ws.Range("N132").Value = -19852.6661

$ws.Range("H136").Value = 2916.0317
$ws.Range("I136").Value = 2500.2554
$ws.Range("J136").Value = 4137.375
$ws.Range("K136").Value = 7500.7662
$ws.Range("L136").Value = 12412.125
$ws.Range("M136").Value = -4950.7662
$ws.Range("N136").Value = -17512.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1540372.9
$ws.Range("I132").Value = 5821
$ws.Range("J132").Value = 4330467.5
$ws.Range("K132").Value = 17463
$ws.Range("L132").Value = 12991402.5
$ws.Range("M132").Value = -14933
$ws.Range("N132").Value = -12996462.5

$ws.Range("H136").Value = 583.6316
$ws.Range("I136").Value = 509.8125
$ws.Range("K136").Value = 1529.4375
$ws.Range("M136").Value = 1020.5625
